# Daily auto push: insert two new rows of data (2026/01/08 and 2026/01/09)
# right before the pre-existing "future placeholder" block, shifting all
# subsequent rows down by two. This mirrors the diff, which turns the
# former row 583 ("2026/12/29 ...") into row 585, etc., and introduces a
# brand-new pair of rows at 583/584.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 583-584; everything from the old row 583 onward
# (previously ending at row 624) shifts down to 585-626.
$ws.Rows("583:584").Insert()

# Row 583: 2026/01/08 (Thu)
$ws.Range("A583").Value = "'2026/01/08"
$ws.Range("A583").Style = "Normal"
$ws.Range("B583").Value = "木"
$ws.Range("C583").Value = 23
$ws.Range("D583").Value = 28

# Row 584: 2026/01/09 (Fri)
$ws.Range("A584").Value = "'2026/01/09"
$ws.Range("A584").Style = "Normal"
$ws.Range("B584").Value = "金"
$ws.Range("C584").Value = 2
$ws.Range("D584").Value = 29
